# Final touches to data mining: refresh the per-city registration stats
# (totals, rates, averages) and a couple of top-brand/model labels that
# shifted after re-running the pipeline on updated source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 186177
$ws.Range("D2").Value = 8.95
$ws.Range("E2").Value = 56.15
$ws.Range("F2").Value = 60.43
$ws.Range("G2").Value = 38.86
$ws.Range("H2").Value = 0.36
$ws.Range("K2").Value = 111.2536865455991
$ws.Range("L2").Value = 161.9785992086307
$ws.Range("M2").Value = 1572.773930184717
$ws.Range("C3").Value = 39470
$ws.Range("D3").Value = 6.56
$ws.Range("E3").Value = 46.33
$ws.Range("F3").Value = 54.62
$ws.Range("G3").Value = 44.63
$ws.Range("H3").Value = 0.26
$ws.Range("K3").Value = 107.1201621484672
$ws.Range("L3").Value = 163.7623576736892
$ws.Range("M3").Value = 1566.245604256397
$ws.Range("C4").Value = 22295
$ws.Range("D4").Value = 7.08
$ws.Range("F4").Value = 54.23
$ws.Range("G4").Value = 45.2
$ws.Range("H4").Value = 0.19
$ws.Range("K4").Value = 106.6654361964566
$ws.Range("L4").Value = 166.1518863302303
$ws.Range("M4").Value = 1561.525140165957
$ws.Range("C5").Value = 18766
$ws.Range("E5").Value = 40.78
$ws.Range("F5").Value = 54.71
$ws.Range("G5").Value = 45.13
$ws.Range("H5").Value = 0.03
$ws.Range("K5").Value = 104.313615048492
$ws.Range("L5").Value = 176.4441618117666
$ws.Range("M5").Value = 1581.429500159864
$ws.Range("C6").Value = 11683
$ws.Range("G6").Value = 58.91
$ws.Range("H6").Value = 0.01
$ws.Range("K6").Value = 101.6444663185825
$ws.Range("L6").Value = 168.7761285574092
$ws.Range("M6").Value = 1571.407943165283
$ws.Range("C7").Value = 7953
$ws.Range("D7").Value = 7.02
$ws.Range("E7").Value = 45.53
$ws.Range("F7").Value = 54.16
$ws.Range("G7").Value = 45.23
$ws.Range("H7").Value = 0.36
$ws.Range("K7").Value = 103.9775179177669
$ws.Range("L7").Value = 163.3787159190853
$ws.Range("M7").Value = 1543.728027159562
$ws.Range("C8").Value = 7096
$ws.Range("D8").Value = 10.01
$ws.Range("E8").Value = 47.29
$ws.Range("F8").Value = 53.13
$ws.Range("G8").Value = 46.39
$ws.Range("H8").Value = 0.13
$ws.Range("K8").Value = 106.1422773393461
$ws.Range("L8").Value = 167.8650748531363
$ws.Range("M8").Value = 1570.619503945885
$ws.Range("C9").Value = 7095
$ws.Range("D9").Value = 8.4
$ws.Range("E9").Value = 43.38
$ws.Range("F9").Value = 53.64
$ws.Range("G9").Value = 45.78
$ws.Range("H9").Value = 0.16
$ws.Range("K9").Value = 104.1938830162086
$ws.Range("L9").Value = 164.5717058024452
$ws.Range("M9").Value = 1549.031994362227
$ws.Range("C10").Value = 6429
$ws.Range("E10").Value = 45
$ws.Range("F10").Value = 58.21
$ws.Range("G10").Value = 40.57
$ws.Range("H10").Value = 0.06
$ws.Range("K10").Value = 102.4209986000933
$ws.Range("L10").Value = 162.55
$ws.Range("M10").Value = 1527.464613470213
$ws.Range("C11").Value = 5603
$ws.Range("D11").Value = 7.44
$ws.Range("E11").Value = 38.85
$ws.Range("F11").Value = 46.23
$ws.Range("G11").Value = 53.26
$ws.Range("H11").Value = 0.27
$ws.Range("K11").Value = 103.1067106907014
$ws.Range("L11").Value = 164.8085163363685
$ws.Range("M11").Value = 1552.078529359272
$ws.Range("P11").Value = 'TOYOTA 7.5'
$ws.Range("C12").Value = 5074
$ws.Range("D12").Value = 7.88
$ws.Range("E12").Value = 39.46
$ws.Range("F12").Value = 51.4
$ws.Range("G12").Value = 48.19
$ws.Range("H12").Value = 0.2
$ws.Range("K12").Value = 103.2046511627907
$ws.Range("L12").Value = 165.0229852440409
$ws.Range("M12").Value = 1539.745368545526
$ws.Range("S12").Value = 'FOCUS 2.6'
$ws.Range("C13").Value = 4618
$ws.Range("D13").Value = 6.8
$ws.Range("E13").Value = 49.78
$ws.Range("F13").Value = 56.37
$ws.Range("G13").Value = 43.27
$ws.Range("H13").Value = 0.17
$ws.Range("K13").Value = 107.0318103074924
$ws.Range("L13").Value = 163.6397287369313
$ws.Range("M13").Value = 1556.445214378519
$ws.Range("C14").Value = 4403
$ws.Range("E14").Value = 37.68
$ws.Range("F14").Value = 48.1
$ws.Range("G14").Value = 51.81
$ws.Range("H14").Value = 0.02
$ws.Range("K14").Value = 102.6298887122417
$ws.Range("L14").Value = 170.4700910273082
$ws.Range("M14").Value = 1567.581648875767
$ws.Range("C15").Value = 2244
$ws.Range("D15").Value = 14.93
$ws.Range("E15").Value = 42.74
$ws.Range("F15").Value = 50.94
$ws.Range("G15").Value = 48.66
$ws.Range("H15").Value = 0.04
$ws.Range("K15").Value = 105.3025846702317
$ws.Range("L15").Value = 175.3071135430916
$ws.Range("M15").Value = 1588.705436720143
$ws.Range("C16").Value = 848
$ws.Range("D16").Value = 11.56
$ws.Range("E16").Value = 39.98
$ws.Range("F16").Value = 56.37
$ws.Range("G16").Value = 43.16
$ws.Range("H16").Value = 0.24
$ws.Range("K16").Value = 103.1770047169811
$ws.Range("L16").Value = 173.2024647887324
$ws.Range("M16").Value = 1547.146226415094
